$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate actual-time tracking for the remaining tasks (rows 12, 14-18)
# Week 1 (F) = 1, Week 2 (G) = 0 for each of these completed tasks.
$rows = @(12, 14, 15, 16, 17, 18)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = 1
    $ws.Cells.Item($r, 7).Value = 0
}

# Move the active selection to D13 (reflects last user interaction point)
$ws.Range("D13").Select()
